# Rename the "CLASSIFICATION" header column to "INFORMATIONDOMAIN".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headerCell = $ws.Range("A1:Z1").Find("CLASSIFICATION")
if ($headerCell) {
    $headerCell.Value = "INFORMATIONDOMAIN"
} else {
    $ws.Range("D1").Value = "INFORMATIONDOMAIN"
}
